# Natmi following Dr Hou advice
# Rebuild the LR-pair result table: the two "Sending cluster" groups (ECs,
# then sCs) are each crossed against all three "Target cluster" values
# (ECs, FAPs, sCs), for a total of 6 result rows (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..D are labels, E..T are the numeric metrics.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs","Fgf16","Fgfr2","ECs",1,0.3333333333333333,0.229066,0.687198,0.1398528162861439,0.1398528162861439,3,1,0.3777296666666667,1.133189,0.2121489992374768,0.2121489992374768,0.08652502382466667,0.7787252144219999,0.02966963501564814,0.02966963501564814),
    @("ECs","Fgf16","Fgfr2","FAPs",1,0.3333333333333333,0.229066,0.687198,0.1398528162861439,0.1398528162861439,3,1,1.092289666666667,3.276869,0.6134761976883921,0.6134761976883921,0.2502064247846666,2.251857823062,0.08579637397123682,0.08579637397123682),
    @("ECs","Fgf16","Fgfr2","sCs",1,0.3333333333333333,0.229066,0.687198,0.1398528162861439,0.1398528162861439,3,1,0.310473,0.931419,0.174374803074131,0.174374803074131,0.071118808218,0.640069273962,0.02438680729925897,0.02438680729925896),
    @("sCs","Fgf16","Fgfr2","ECs",3,1,1.408841666666667,4.226525,0.860147183713856,0.8601471837138561,3,1,0.3777296666666667,1.133189,0.2121489992374768,0.2121489992374768,0.532161293136111,4.789451638225,0.1824793642218287,0.1824793642218287),
    @("sCs","Fgf16","Fgfr2","FAPs",3,1,1.408841666666667,4.226525,0.860147183713856,0.8601471837138561,3,1,1.092289666666667,3.276869,0.6134761976883921,0.6134761976883921,1.538863194469444,13.849768750225,0.5276798237171553,0.5276798237171554),
    @("sCs","Fgf16","Fgfr2","sCs",3,1,1.408841666666667,4.226525,0.860147183713856,0.8601471837138561,3,1,0.310473,0.931419,0.174374803074131,0.174374803074131,0.437407298775,3.936665688975,0.149987995774872,0.149987995774872)
)

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $r).Value = $row[$i]
    }
    $r++
}
